# Update column F ("dSF") values on Sheet1 to match the repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 3
    3  = -1
    4  = -1
    5  = -3
    6  = -1
    7  = 2
    8  = -1
    9  = -2
    10 = -1
    11 = 0
    12 = -1
    13 = -3
    14 = 1
    15 = -3
    16 = -3
    17 = 3
    18 = -2
    19 = 4
    20 = 1
    21 = 1
    23 = 3
    24 = -1
    25 = 5
    26 = 7
    29 = -1
    30 = 3
    32 = 5
    33 = -2
    34 = 5
    35 = 3
    37 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
